# Add "Magnetometer sensor data" (columns I/J/K) to the sensor data worksheet,
# mirroring the existing "raw front sensor data" (A/B/C) and
# "raw side sensor data" (E/F/G) blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Raw magnetometer readings (column I, rows 2-50) ---
$magValues = @(358,360,361,358,353,358,359,357,355,355,359,359,359,360,357,358,352,359,359,355,359,358,357,361,353,358,355,358,359,359,354,355,357,353,357,360,359,361,358,357,359,352,360,358,358,354,360,357,355)

for ($i = 0; $i -lt $magValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $magValues[$i]
}

# --- Header / label / summary block (column I/J/K, rows 1-3) ---
$ws.Range("I1").Value = "Magnetometer sensor data"
$ws.Range("J1").Value = "actual"
$ws.Range("K1").Value = 360

$ws.Range("J2").Value = "std dev"
$ws.Range("K2").Formula = "=STDEV.P(I2:I50)"

$ws.Range("J3").Value = "mean"
$ws.Range("K3").Formula = "=AVERAGE(I2:I50)"

# --- Right-align the "label" columns (B, F, J) to match existing B/F styling ---
$ws.Range("B1").HorizontalAlignment = -4152
$ws.Range("B2").HorizontalAlignment = -4152
$ws.Range("B3").HorizontalAlignment = -4152
$ws.Range("F1").HorizontalAlignment = -4152
$ws.Range("F2").HorizontalAlignment = -4152
$ws.Range("F3").HorizontalAlignment = -4152
$ws.Range("J1").HorizontalAlignment = -4152
$ws.Range("J2").HorizontalAlignment = -4152
$ws.Range("J3").HorizontalAlignment = -4152

# --- Column widths (B, F, I, J) ---
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(6).ColumnWidth = 14.166666666666666
$ws.Columns.Item(9).ColumnWidth = 25.498697916666668
$ws.Columns.Item(10).ColumnWidth = 9.998697916666666

# --- Match the saved selection / active cell ---
[void]$ws.Range("K4").Select()
